$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.447.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.85"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.30"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3643"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.249"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.640"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001255"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.290"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.634.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.13"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06941"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.549"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.459.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.244"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.455"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.76"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.306"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.97"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.318"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.818.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.895"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("E34").Value = "  +7.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9657"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02861"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.275"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2558"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07280"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08876"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7126"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6564"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.358"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08001"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.219"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.37%  "
